$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 349, shifting the existing rows 349:374 down to 350:375
$ws.Rows.Item(349).Insert()

# Populate the newly inserted row 349 with this week's new data point
$ws.Range("A349").Value = 9
$ws.Range("B349").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C349").Value = "Metropolitana"
$ws.Range("D349").Value = 45106
$ws.Range("E349").Value = 13
$ws.Range("F349").Value = 100112001
$ws.Range("G349").Value = "Berenjena"
$ws.Range("H349").Value = "Sin especificar"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 70
$ws.Range("K349").Value = 7000
$ws.Range("L349").Value = 8000
$ws.Range("M349").Value = 7500
$ws.Range("N349").Value = "$/caja 50 unidades"
$ws.Range("O349").Value = "Región de Arica y Parinacota"
$ws.Range("P349").Value = 150
$ws.Range("Q349").Value = 50
$ws.Range("R349").Value = "Hortaliza"
